$wb = $excel.ActiveWorkbook

# Rename the 4th sheet (味全龍) to "Dragons"
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "Dragons"

# Make the Dragons sheet the active/selected tab
$ws4.Activate()
$ws4.Select()
